# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff"
# - Refresh the "Latest Handoff/HO Xliff Generate Date" timestamps
# - Widen the Status columns so the new, longer label fits

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status column: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 (Latest HO Xliff Generate Date) and de-de!H2 (Latest Handoff Datetime)
# shared the same "2016-09-04 17:03:58" value -> bump to 17:04:32
$overview.Range("G2").Value = "2016-09-04 17:04:32"
$dede.Range("H2").Value     = "2016-09-04 17:04:32"

# zh-cn!H2 (Latest Handoff Datetime) had its own "2016-09-04 17:03:54" -> bump to 17:04:28
$zhcn.Range("H2").Value = "2016-09-04 17:04:28"

# --- Column widths for the Status columns (wider to fit "Ready for handoff") ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth     = 16.33
$dede.Columns.Item(3).ColumnWidth     = 16.33
